$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15, Col C: "BICC" -> "Arquitectura"
$ws.Range("C15").Value = "Arquitectura"

# Row 21, Col C: clear the cell (remove "Pool Visualizacion")
$ws.Range("C21").Value = ""

# Row 23, Col C: "MKT Mad" -> "MKT"
$ws.Range("C23").Value = "MKT"

# Row 24, Col C: "MKT BIO" -> "MKT"
$ws.Range("C24").Value = "MKT"

# Row 25, Col C: "MKT DnA" -> "MKT"
$ws.Range("C25").Value = "MKT"

# Row 14: apply new font style (Times New Roman 12) and row height
$ws.Range("A14").Font.Name = "Times New Roman"
$ws.Range("A14").Font.Size = 12
$ws.Rows("14").RowHeight = 15.5

# Sheet view changes
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("C15").Select()
